# Ajout de spécifications dans les scénarios
# Applies the target edit to ScenariosGuest.xlsx:
#  - Replaces the generic "Oui -->" / "Non -->" labels with more specific
#    outcome labels on the "Création d'un compte", "Connexion" and
#    "Filtrage des locations" sheets, and adds the corresponding extra
#    rows/reactions.
#  - Updates the selected cell on several sheets and moves the active tab
#    from "Consultation des locations" to "Filtrage des locations".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Création d'un compte"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Création d'un compte")

$ws1.Range("B8").Value  = "Toutes les conditions sont OK"
$ws1.Range("C9").Value  = "Charge la page d'accueil"
$ws1.Range("B10").Value = ""
$ws1.Range("C10").Value = "Connecte l'utilisateur au compte créé"
$ws1.Range("B11").Value = "Les conditions ne sont pas toutes OK"
$ws1.Range("C11").Value = "Recharge la page"
$ws1.Range("C12").Value = "Affiche un message d'erreur"

$ws1.Range("B35").Select()

# ---------------------------------------------------------------------
# Sheet 2: "Connexion"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Connexion")

$ws2.Range("B6").Value  = "Toutes les conditions sont OK"
$ws2.Range("B7").Value  = ""
$ws2.Range("C7").Value  = "Charge la page d'acceuil"
$ws2.Range("B8").Value  = "Les conditions ne sont pas toutes OK"
$ws2.Range("C8").Value  = "Recharge la page"
$ws2.Range("C9").Value  = "Affiche un message d'erreur"

$ws2.Range("C10").Select()

# ---------------------------------------------------------------------
# Sheet 4: "Consultation des locations" (loses tabSelected)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Consultation des locations")
$ws4.Range("C19").Select()

# ---------------------------------------------------------------------
# Sheet 6: "Filtrage des locations" (becomes the active tab)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Filtrage des locations")

$ws6.Range("C4").Value = "Affiche les locations présentes correspondant au filtrage"
$ws6.Range("B5").Value = ""
$ws6.Range("C5").Value = ""
$ws6.Range("B6").Value = "Les valeurs des champs de filtrage ne sont pas toutes valides"

# Row 5 no longer needs the taller "Oui/Non" wrapped-text height, row 6 now
# needs it instead (its text got longer).
$ws6.Rows.Item(5).AutoFit()
$ws6.Rows.Item(6).RowHeight = 30

$ws6.Activate()
$ws6.Range("B7").Select()
